# Reorders data rows 2-25 (columns A:F) of the active sheet according to a
# fixed permutation determined by the target workbook state.
# Row 1 (header) and row 26 (totals) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of destination row -> source row (both are original, pre-edit row numbers)
$map = @{
    2  = 9
    3  = 14
    4  = 2
    5  = 8
    6  = 5
    7  = 13
    8  = 7
    9  = 15
    10 = 10
    11 = 11
    12 = 3
    13 = 4
    14 = 12
    15 = 6
    16 = 18
    17 = 21
    18 = 19
    19 = 20
    20 = 17
    21 = 16
    22 = 23
    23 = 22
    24 = 24
    25 = 25
}

# Snapshot all original values for rows 2-25, columns A-F (1-6), before
# writing anything, since several rows both read-from and are written-to.
$snapshot = @{}
for ($r = 2; $r -le 25; $r++) {
    $rowVals = @()
    for ($c = 1; $c -le 6; $c++) {
        $rowVals += , ($ws.Cells.Item($r, $c).Value2)
    }
    $snapshot[$r] = $rowVals
}

# Write back the values according to the mapping.
foreach ($destRow in $map.Keys) {
    $srcRow = $map[$destRow]
    $vals = $snapshot[$srcRow]
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($destRow, $c).Value2 = $vals[$c - 1]
    }
}
